$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 1.05
$ws.Range("J4").Value = 12
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 34
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.95
$ws.Range("W4").Value = 51
$ws.Range("X4").Value = 101
$ws.Range("AA4").Value = 126
$ws.Range("AD4").Value = 21
$ws.Range("AE4").Value = 29
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 201
$ws.Range("AI4").Value = 11
$ws.Range("AJ4").Value = 13
$ws.Range("AL4").Value = 11
$ws.Range("AM4").Value = 26
$ws.Range("AN4").Value = 17
$ws.Range("AO4").Value = 51
$ws.Range("AQ4").Value = 301
$ws.Range("AS4").Value = 201
$ws.Range("AW4").Value = 351
$ws.Range("AY4").Value = 4.75
$ws.Range("BA4").Value = 8

# Row 5
$ws.Range("I5").Value = 3.1
$ws.Range("Q5").Value = 1.75
$ws.Range("R5").Value = 2.05
$ws.Range("X5").Value = 12
$ws.Range("AO5").Value = 12
$ws.Range("AX5").Value = 5

# Row 6
$ws.Range("H6").Value = 4.2
$ws.Range("J6").Value = 2.2
$ws.Range("K6").Value = 2.37
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.7
$ws.Range("R6").Value = 2.1
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("U6").Value = 1.75
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 8
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 7.5
$ws.Range("AH6").Value = 15
$ws.Range("AO6").Value = 8
$ws.Range("AT6").Value = 3.25
$ws.Range("BC6").Value = 201

# Row 7
$ws.Range("G7").Value = 3.2
$ws.Range("H7").Value = 3.7
$ws.Range("I7").Value = 2.05
$ws.Range("J7").Value = 3.75
$ws.Range("L7").Value = 2.75
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 3.75
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 2.05
$ws.Range("W7").Value = 11
$ws.Range("X7").Value = 17
$ws.Range("Y7").Value = 11
$ws.Range("Z7").Value = 34
$ws.Range("AA7").Value = 23
$ws.Range("AI7").Value = 11
$ws.Range("AJ7").Value = 9
$ws.Range("AO7").Value = 17
$ws.Range("BA7").Value = 41

# Row 8
$ws.Range("G8").Value = 7.5
$ws.Range("H8").Value = 4.75
$ws.Range("I8").Value = 1.38
$ws.Range("J8").Value = 6
$ws.Range("K8").Value = 2.75
$ws.Range("L8").Value = 1.83
$ws.Range("N8").Value = 21
$ws.Range("O8").Value = 1.11
$ws.Range("P8").Value = 6.5
$ws.Range("Q8").Value = 1.4
$ws.Range("R8").Value = 2.88
$ws.Range("S8").Value = 1.22
$ws.Range("T8").Value = 4
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.2
$ws.Range("W8").Value = 23
$ws.Range("X8").Value = 41
$ws.Range("Y8").Value = 21
$ws.Range("Z8").Value = 81
$ws.Range("AC8").Value = 21
$ws.Range("AD8").Value = 10
$ws.Range("AE8").Value = 15
$ws.Range("AG8").Value = 151
$ws.Range("AI8").Value = 8.5
$ws.Range("AJ8").Value = 9
$ws.Range("AK8").Value = 10
$ws.Range("AL8").Value = 10
$ws.Range("AM8").Value = 21
$ws.Range("AN8").Value = 8.5
$ws.Range("AO8").Value = 34
$ws.Range("AP8").Value = 29
$ws.Range("AQ8").Value = 101
$ws.Range("AR8").Value = 101
$ws.Range("AT8").Value = 4
$ws.Range("AU8").Value = 8
$ws.Range("AY8").Value = 6.5
$ws.Range("AZ8").Value = 13
$ws.Range("BA8").Value = 17

# Row 9
$ws.Range("G9").Value = 1.25
$ws.Range("H9").Value = 5.5
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 1.73
$ws.Range("K9").Value = 2.63
$ws.Range("L9").Value = 8.5
$ws.Range("N9").Value = 17
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.6
$ws.Range("R9").Value = 2.3
$ws.Range("S9").Value = 1.29
$ws.Range("T9").Value = 3.5
$ws.Range("U9").Value = 2.05
$ws.Range("V9").Value = 1.7
$ws.Range("W9").Value = 7.5
$ws.Range("X9").Value = 6.5
$ws.Range("Y9").Value = 9
$ws.Range("Z9").Value = 7.5
$ws.Range("AB9").Value = 29
$ws.Range("AD9").Value = 11
$ws.Range("AH9").Value = 26
$ws.Range("AI9").Value = 51
$ws.Range("AJ9").Value = 29
$ws.Range("AK9").Value = 151
$ws.Range("AL9").Value = 67
$ws.Range("AM9").Value = 67
$ws.Range("AO9").Value = 5.5
$ws.Range("AQ9").Value = 15
$ws.Range("AR9").Value = 41
$ws.Range("AS9").Value = 126
$ws.Range("AT9").Value = 3.5
$ws.Range("AU9").Value = 9.5
$ws.Range("AV9").Value = 51
$ws.Range("AX9").Value = 10
$ws.Range("AY9").Value = 41
$ws.Range("AZ9").Value = 41
$ws.Range("BA9").Value = 201
$ws.Range("BB9").Value = 201
$ws.Range("BC9").Value = 351

# Row 10
$ws.Range("G10").Value = 1.75
$ws.Range("H10").Value = 3.9
$ws.Range("I10").Value = 4.33
$ws.Range("J10").Value = 2.3
$ws.Range("K10").Value = 2.4
$ws.Range("L10").Value = 4.33
$ws.Range("Q10").Value = 1.53
$ws.Range("R10").Value = 2.4
$ws.Range("U10").Value = 1.57
$ws.Range("V10").Value = 2.25
$ws.Range("X10").Value = 10
$ws.Range("Z10").Value = 15
$ws.Range("AA10").Value = 13
$ws.Range("AB10").Value = 19
$ws.Range("AC10").Value = 17
$ws.Range("AD10").Value = 8
$ws.Range("AE10").Value = 13
$ws.Range("AF10").Value = 41
$ws.Range("AH10").Value = 17
$ws.Range("AI10").Value = 23
$ws.Range("AL10").Value = 29
$ws.Range("AN10").Value = 4
$ws.Range("AO10").Value = 9
$ws.Range("AP10").Value = 15
$ws.Range("AQ10").Value = 26
$ws.Range("AS10").Value = 81
$ws.Range("AU10").Value = 7.5
$ws.Range("AX10").Value = 6.5
$ws.Range("AY10").Value = 21
$ws.Range("AZ10").Value = 23
$ws.Range("BA10").Value = 67

# Row 25
$ws.Range("G25").Value = 3.25
$ws.Range("I25").Value = 2.15
$ws.Range("J25").Value = 3.6
$ws.Range("L25").Value = 2.75
$ws.Range("S25").Value = 1.3
$ws.Range("T25").Value = 3.4
$ws.Range("U25").Value = 1.53
$ws.Range("V25").Value = 2.38
$ws.Range("AE25").Value = 11
$ws.Range("AF25").Value = 34
$ws.Range("AH25").Value = 10
$ws.Range("AI25").Value = 12
$ws.Range("AK25").Value = 21
$ws.Range("AP25").Value = 21
$ws.Range("AT25").Value = 3.4
$ws.Range("AU25").Value = 7
$ws.Range("AW25").Value = 401
$ws.Range("BA25").Value = 41
